# The footer's "... of <NUMPAGES>" field was stored as a <w:fldSimple> shortcut
# (a single-run "simple field"). Re-express it using the general begin/separate/end
# <w:fldChar> run sequence (the same form already used for the "Page <PAGE>" field
# earlier in the same paragraph), fixing the inconsistent field-XML formatting.
$d = $word.ActiveDocument

$sec = $d.Sections(1)

$footer = $sec.Footers(1)  # wdHeaderFooterPrimary
$footerXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Footer"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Page </w:t>
  </w:r>
  <w:r>
    <w:fldChar w:fldCharType="begin"/>
  </w:r>
  <w:r>
    <w:instrText xml:space="preserve"> PAGE   \* MERGEFORMAT </w:instrText>
  </w:r>
  <w:r>
    <w:fldChar w:fldCharType="separate"/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:noProof/>
    </w:rPr>
    <w:t>1</w:t>
  </w:r>
  <w:r>
    <w:fldChar w:fldCharType="end"/>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> of </w:t>
  </w:r>
  <w:r>
    <w:fldChar w:fldCharType="begin"/>
  </w:r>
  <w:r>
    <w:instrText xml:space="preserve"> NUMPAGES   \* MERGEFORMAT </w:instrText>
  </w:r>
  <w:r>
    <w:fldChar w:fldCharType="separate"/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:noProof/>
    </w:rPr>
    <w:t>4</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:noProof/>
    </w:rPr>
    <w:fldChar w:fldCharType="end"/>
  </w:r>
</w:p>
'@
$footer.Range.InsertXML($footerXml) | Out-Null

# The first-page header's paragraph contained a stray empty run (<w:r><w:t/></w:r>)
# that carried no actual text. Drop it, leaving just the empty, styled paragraph.
$firstPageHeader = $sec.Headers(2)  # wdHeaderFooterFirstPage
$headerXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Header"/>
  </w:pPr>
</w:p>
'@
$firstPageHeader.Range.InsertXML($headerXml) | Out-Null
